# Updates cryptos list price / volume(1h) data, including a couple of
# row re-orderings (Hedera/Kaspa and TheGraph/Maker swapped places).
# Values that would otherwise be auto-parsed by Excel as numbers (and so
# lose their original text formatting, e.g. trailing zeros) are entered
# with a leading apostrophe to force a text/quoted-text cell, matching
# the original inline-string (t="inlineStr") cell content exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.786.33"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.929.93"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'375.94"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").Value = "'100.05"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'35.81"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "'0.0848"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "3.390.78"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "'17.99"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'11.78"
$ws.Range("E16").Value = "  +59.42%  "
$ws.Range("D17").Value = "2.927.99"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "'0.991"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "50.748.76"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "'3.03"
$ws.Range("E20").Value = "  -7.40%  "
$ws.Range("D21").Value = "'12.34"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "'69.16"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'265.90"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  +7.03%  "
$ws.Range("D26").Value = "'7.97"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").Value = "'7.22"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'25.45"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.162"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.108"
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "'50.65"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").Value = "'32.84"
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("D36").Value = "'0.0429"
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'16.30"
$ws.Range("E40").Value = "  -5.11%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "'2.44"
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("D43").Value = "'119.28"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'20.99"
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("D45").Value = "'3.37"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.993.42"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.260"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("D50").Value = "'0.0320"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("D51").Value = "'5.22"
$ws.Range("E51").Value = "  +2.11%  "
